$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("I2").Value = 4
$ws.Range("J2").Value = 0.0002777777777777778
$ws.Range("K2").Value = 1922
$ws.Range("L2").Value = 0.003844
